$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(121).Insert()

$ws.Cells.Item(121, 1).Value = 4
$ws.Cells.Item(121, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(121, 3).Value = "Los Lagos"
$ws.Cells.Item(121, 4).Value = 44589
$ws.Cells.Item(121, 5).Value = 10
$ws.Cells.Item(121, 6).Value = 100112037
$ws.Cells.Item(121, 7).Value = "Cebollín"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 160
$ws.Cells.Item(121, 11).Value = 6000
$ws.Cells.Item(121, 12).Value = 6500
$ws.Cells.Item(121, 13).Value = 6250
$ws.Cells.Item(121, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(121, 15).Value = "Región Metropolitana"
$ws.Cells.Item(121, 16).Value = 174
$ws.Cells.Item(121, 17).Value = 36
$ws.Cells.Item(121, 18).Value = "Hortaliza"
